$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.2211538461538461
$ws.Range("C2").Value = 0.5192307692307693
$ws.Range("J2").Value = 0.02163461538461538
$ws.Range("P2").Value = 0.1490384615384615
$ws.Range("S2").Value = 0.0889423076923077
$ws.Range("B3").Value = 0.009009009009009009
$ws.Range("C3").Value = 0.02702702702702703
$ws.Range("J3").Value = 0.009009009009009009
$ws.Range("O3").Value = 0.004504504504504504
$ws.Range("P3").Value = 0.7837837837837838
$ws.Range("S3").Value = 0.1666666666666667
$ws.Range("J4").Value = 0.0425531914893617
$ws.Range("P4").Value = 0.7659574468085106
$ws.Range("S4").Value = 0.1914893617021277
$ws.Range("B6").Value = 0.0811965811965812
$ws.Range("D6").Value = 0.01282051282051282
$ws.Range("E6").Value = 0.004273504273504274
$ws.Range("F6").Value = 0.09829059829059829
$ws.Range("J6").Value = 0.217948717948718
$ws.Range("O6").Value = 0.02991452991452992
$ws.Range("Q6").Value = 0.1837606837606838
$ws.Range("R6").Value = 0.07264957264957266
$ws.Range("S6").Value = 0.2991452991452991
$ws.Range("B7").Value = 0.06532663316582915
$ws.Range("D7").Value = 0.01005025125628141
$ws.Range("F7").Value = 0.03517587939698492
$ws.Range("J7").Value = 0.1206030150753769
$ws.Range("O7").Value = 0.03517587939698492
$ws.Range("Q7").Value = 0.185929648241206
$ws.Range("R7").Value = 0.1005025125628141
$ws.Range("S7").Value = 0.4472361809045226
$ws.Range("B8").Value = 0.1184210526315789
$ws.Range("D8").Value = 0.02412280701754386
$ws.Range("F8").Value = 0.07675438596491228
$ws.Range("J8").Value = 0.1140350877192982
$ws.Range("O8").Value = 0.01754385964912281
$ws.Range("Q8").Value = 0.1907894736842105
$ws.Range("R8").Value = 0.05482456140350877
$ws.Range("S8").Value = 0.4035087719298245
$ws.Range("B9").Value = 0.1291512915129151
$ws.Range("D9").Value = 0.003690036900369004
$ws.Range("F9").Value = 0.06642066420664207
$ws.Range("J9").Value = 0.09225092250922509
$ws.Range("O9").Value = 0.02583025830258303
$ws.Range("Q9").Value = 0.1697416974169742
$ws.Range("R9").Value = 0.05904059040590406
$ws.Range("S9").Value = 0.4538745387453875
$ws.Range("B10").Value = 0.1307490144546649
$ws.Range("D10").Value = 0.02102496714848883
$ws.Range("E10").Value = 0.001314060446780552
$ws.Range("F10").Value = 0.05387647831800263
$ws.Range("J10").Value = 0.1268068331143233
$ws.Range("O10").Value = 0.01576872536136662
$ws.Range("Q10").Value = 0.2345597897503285
$ws.Range("R10").Value = 0.06701708278580815
$ws.Range("S10").Value = 0.3488830486202366
$ws.Range("G11").Value = 0.1477987421383648
$ws.Range("J11").Value = 0.1006289308176101
$ws.Range("K11").Value = 0.2232704402515723
$ws.Range("L11").Value = 0.5157232704402516
$ws.Range("S11").Value = 0.01257861635220126
$ws.Range("G12").Value = 0.7341040462427746
$ws.Range("J12").Value = 0.1849710982658959
$ws.Range("K12").Value = 0.01734104046242774
$ws.Range("L12").Value = 0.04624277456647399
$ws.Range("S12").Value = 0.01734104046242774
$ws.Range("F13").Value = 0.025
$ws.Range("G13").Value = 0.725
$ws.Range("J13").Value = 0.15
$ws.Range("S13").Value = 0.1
$ws.Range("H15").Value = 0.1173708920187793
$ws.Range("I15").Value = 0.04225352112676056
$ws.Range("J15").Value = 0.3802816901408451
$ws.Range("K15").Value = 0.09389671361502347
$ws.Range("M15").Value = 0.004694835680751174
$ws.Range("O15").Value = 0.05633802816901409
$ws.Range("S15").Value = 0.3051643192488263
$ws.Range("F16").Value = 0.02247191011235955
$ws.Range("H16").Value = 0.1235955056179775
$ws.Range("I16").Value = 0.1198501872659176
$ws.Range("J16").Value = 0.4269662921348314
$ws.Range("K16").Value = 0.08614232209737828
$ws.Range("M16").Value = 0.01872659176029963
$ws.Range("O16").Value = 0.0299625468164794
$ws.Range("S16").Value = 0.1722846441947566
$ws.Range("F17").Value = 0.02120141342756184
$ws.Range("H17").Value = 0.1360424028268551
$ws.Range("I17").Value = 0.1325088339222615
$ws.Range("J17").Value = 0.4858657243816255
$ws.Range("K17").Value = 0.07597173144876325
$ws.Range("M17").Value = 0.008833922261484099
$ws.Range("N17").Value = 0.00176678445229682
$ws.Range("O17").Value = 0.04063604240282685
$ws.Range("S17").Value = 0.09717314487632508
$ws.Range("F18").Value = 0.005649717514124294
$ws.Range("H18").Value = 0.1355932203389831
$ws.Range("I18").Value = 0.1186440677966102
$ws.Range("J18").Value = 0.4915254237288136
$ws.Range("K18").Value = 0.1016949152542373
$ws.Range("M18").Value = 0.01694915254237288
$ws.Range("O18").Value = 0.05084745762711865
$ws.Range("S18").Value = 0.07909604519774012
$ws.Range("F19").Value = 0.01144492131616595
$ws.Range("H19").Value = 0.2160228898426323
$ws.Range("I19").Value = 0.09585121602288985
$ws.Range("J19").Value = 0.3969957081545064
$ws.Range("K19").Value = 0.09585121602288985
$ws.Range("M19").Value = 0.01859799713876967
$ws.Range("N19").Value = 0.000715307582260372
$ws.Range("O19").Value = 0.05579399141630902
$ws.Range("S19").Value = 0.1087267525035765
